$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text representation instead of
# being auto-converted to a number by Excel (e.g. "1.000" -> 1).
$ws.Range("D2:D51").NumberFormat = "@"

# --- Row 21 / 22 swap: Uniswap and WrappedliquidstakedEther2.0 exchange rows ---
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.535"
$ws.Range("E21").Value = "  +6.54%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.168.85"
$ws.Range("E22").Value = "  +0.10%  "

# --- Price / Volume updates for remaining rows ---
$ws.Range("D2").Value = "30.587.25"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.920.46"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "247.18"
$ws.Range("E5").Value = "  +2.83%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "0.4738"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "0.2887"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").Value = "0.06814"
$ws.Range("E9").Value = "  +3.56%  "
$ws.Range("D10").Value = "104.87"
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").Value = "18.36"
$ws.Range("E11").Value = "  -4.76%  "
$ws.Range("D12").Value = "1.915.89"
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D13").Value = "0.07705"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "5.276"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "0.6686"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "291.69"
$ws.Range("E16").Value = "  -3.76%  "
$ws.Range("D17").Value = "30.583.37"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "0.000007593"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "12.93"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "6.368"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "9.385"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").Value = "167.83"
$ws.Range("E26").Value = "  +1.77%  "
$ws.Range("D27").Value = "21.05"
$ws.Range("E27").Value = "  +6.24%  "
$ws.Range("D28").Value = "2.116"
$ws.Range("E28").Value = "  +3.94%  "
$ws.Range("D29").Value = "0.1065"
$ws.Range("E29").Value = "  -5.28%  "
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("D31").Value = "4.169"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "4.058"
$ws.Range("E32").Value = "  +3.18%  "
$ws.Range("D33").Value = "0.05028"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").Value = "1.143"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("D36").Value = "0.02065"
$ws.Range("E36").Value = "  +4.91%  "
$ws.Range("D37").Value = "2.744"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "2.690"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "2.049"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("D40").Value = "110.86"
$ws.Range("E40").Value = "  +3.60%  "
$ws.Range("D41").Value = "0.8750"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "0.4383"
$ws.Range("E42").Value = "  +5.59%  "
$ws.Range("D43").Value = "5.854"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "67.24"
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("D46").Value = "7.238"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "9.328"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "48.21"
$ws.Range("E48").Value = "  +13.87%  "
$ws.Range("D49").Value = "0.1231"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("D50").Value = "34.77"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  +4.88%  "
